$wb = $excel.ActiveWorkbook

# Fix the typo in the second worksheet's name.
$ws = $wb.Worksheets.Item(2)
$ws.Name = "시료채취지점"

# Make that sheet the active one (it was previously the last sheet).
$ws.Activate()
